# Update countries & provincias Spain
# Refresh the COVID-19 snapshot data: update the "last updated" timestamp,
# swap the ranking of Burkina Faso/Uruguay and Timor Oriental/Santa Lucia
# (their case counts crossed over), and refresh the numeric columns for the
# rows whose underlying source data changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Septiembre de 2020 a las 02:30"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 7320669
$ws.Range("C4").Value = 33108
$ws.Range("D4").Value = 4559998
$ws.Range("E4").Value = 2551218
$ws.Range("G4").Value = 276
$ws.Range("H4").Value = 209453

# --- Brasil (row 6) ---
$ws.Range("E6").Value = 530445
$ws.Range("G6").Value = 335
$ws.Range("H6").Value = 141776

# --- Peru (row 9) ---
$ws.Range("B9").Value = 805302
$ws.Range("C9").Value = 5160
$ws.Range("D9").Value = 664490
$ws.Range("E9").Value = 108550
$ws.Range("G9").Value = 120
$ws.Range("H9").Value = 32262

# --- Chequia (row 57) ---
$ws.Range("B57").Value = 64597
$ws.Range("C57").Value = 1303
$ws.Range("D57").Value = 31268
$ws.Range("E57").Value = 32723
$ws.Range("G57").Value = 15
$ws.Range("H57").Value = 606

# --- Noruega (row 93) ---
$ws.Range("B93").Value = 13698
$ws.Range("C93").Value = 71
$ws.Range("E93").Value = 2238

# --- Congo (row 126) ---
$ws.Range("B126").Value = 5008
$ws.Range("C126").Value = 3
$ws.Range("E126").Value = 1032

# --- Surinam (row 127) ---
$ws.Range("B127").Value = 4835
$ws.Range("C127").Value = 4
$ws.Range("D127").Value = 4661
$ws.Range("E127").Value = 72

# --- Benin (row 150) ---
$ws.Range("B150").Value = 2340
$ws.Range("C150").Value = 15
$ws.Range("E150").Value = 340

# --- Rows 154/155: Uruguay and Burkina Faso swap rank, with refreshed data ---
$ws.Range("A154").Value = "Burkina Faso"
$ws.Range("C154").Value = 35
$ws.Range("D154").Value = 1276
$ws.Range("E154").Value = 676
$ws.Range("H154").Value = 56

$ws.Range("A155").Value = "Uruguay"
$ws.Range("B155").Value = 2008
$ws.Range("C155").Value = 10
$ws.Range("D155").Value = 1728
$ws.Range("E155").Value = 233
$ws.Range("H155").Value = 47

# --- Togo (row 159) ---
$ws.Range("B159").Value = 1743
$ws.Range("C159").Value = 7
$ws.Range("D159").Value = 1330
$ws.Range("E159").Value = 367

# --- Niger (row 166) ---
$ws.Range("B166").Value = 1196
$ws.Range("C166").Value = 2
$ws.Range("E166").Value = 20

# --- Rows 207/208: Santa Lucia and Timor Oriental swap rank (data unchanged) ---
$ws.Range("A207").Value = "Timor Oriental"
$ws.Range("A208").Value = "Santa Lucia"
